$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "DO_mg_l"
$ws.Range("B6").Value = "TDS_mg_l"
$ws.Range("B7").Value = "Salinity_ppt"

$wb.Save()
